$d = $word.ActiveDocument

$replacements = @(
    @("14×61=", "92×24="),
    @("38×79=", "92×98="),
    @("82×53=", "86×23="),
    @("36×59=", "34×55="),
    @("55×36=", "85×92="),
    @("52×44=", "99×96="),
    @("66×85=", "22×57="),
    @("98×41=", "80×46="),
    @("43×50=", "34×61="),
    @("96×80=", "28×96="),
    @("71×45=", "51×37="),
    @("94×21=", "87×25="),
    @("37×37=", "72×59="),
    @("55×28=", "19×65="),
    @("29×67=", "46×36="),
    @("98×53=", "64×99="),
    @("44×46=", "12×87="),
    @("47×38=", "65×76="),
    @("17×86=", "34×64="),
    @("80×85=", "54×60="),
    @("78×94=", "29×16="),
    @("44×14=", "42×79="),
    @("66×63=", "77×20="),
    @("98×18=", "80×35="),
    @("29×27=", "69×89=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
